$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2177.3333
$ws.Range("I19").Value = 2652.0715
$ws.Range("J19").Value = 1666.0769
$ws.Range("K19").Value = 2652.0715
$ws.Range("L19").Value = 1666.0769
$ws.Range("M19").Value = -2477.0715
$ws.Range("N19").Value = -2016.0769
$ws.Range("H53").Value = 216.78947
$ws.Range("I53").Value = 187
$ws.Range("K53").Value = 187
$ws.Range("M53").Value = 450
$ws.Range("H62").Value = 11215.414
$ws.Range("I62").Value = 8434.105
$ws.Range("J62").Value = 16499.9
$ws.Range("K62").Value = 8434.105
$ws.Range("L62").Value = 16499.9
$ws.Range("M62").Value = -7810.105
$ws.Range("N62").Value = -17747.9
$ws.Range("H64").Value = 6154.375
$ws.Range("I64").Value = 5608.75
$ws.Range("J64").Value = 6700
$ws.Range("K64").Value = 5608.75
$ws.Range("L64").Value = 6700
$ws.Range("M64").Value = -5360.75
$ws.Range("N64").Value = -7196
$ws.Range("H65").Value = 11215.414
$ws.Range("I65").Value = 8434.105
$ws.Range("J65").Value = 16499.9
$ws.Range("K65").Value = 42170.52499999999
$ws.Range("L65").Value = 82499.5
$ws.Range("M65").Value = -39050.52499999999
$ws.Range("N65").Value = -88739.5
$ws.Range("H67").Value = 6154.375
$ws.Range("I67").Value = 5608.75
$ws.Range("J67").Value = 6700
$ws.Range("K67").Value = 5608.75
$ws.Range("L67").Value = 6700
$ws.Range("M67").Value = -4750.75
$ws.Range("N67").Value = -8416
$ws.Range("H117").Value = 80742
$ws.Range("J117").Value = 80742
$ws.Range("L117").Value = 80742
$ws.Range("N117").Value = -89920
$ws.Range("H127").Value = 2443.1428
$ws.Range("I127").Value = 2647.5
$ws.Range("J127").Value = 1217
$ws.Range("K127").Value = 7942.5
$ws.Range("L127").Value = 3651
$ws.Range("M127").Value = -2982.5
$ws.Range("N127").Value = -13571
$ws.Range("H132").Value = 104540.56
$ws.Range("I132").Value = 113306.09
$ws.Range("K132").Value = 339918.27
$ws.Range("M132").Value = -337388.27
$ws.Range("H137").Value = 1017645.1
$ws.Range("I137").Value = 3100.8572
$ws.Range("K137").Value = 9302.571599999999
$ws.Range("M137").Value = -6752.571599999999
$ws.Range("H138").Value = 2178.4856
$ws.Range("I138").Value = 1322.2667
$ws.Range("J138").Value = 2820.65
$ws.Range("K138").Value = 3966.800099999999
$ws.Range("L138").Value = 8461.950000000001
$ws.Range("M138").Value = 1173.199900000001
$ws.Range("N138").Value = -18741.95

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4167385
$ws.Range("I61").Value = 4167385
$ws.Range("K61").Value = 4167385
$ws.Range("M61").Value = -4167173
$ws.Range("H110").Value = 2155.8215
$ws.Range("I110").Value = 1773.875
$ws.Range("J110").Value = 4447.5
$ws.Range("K110").Value = 1773.875
$ws.Range("L110").Value = 4447.5
$ws.Range("M110").Value = 271.125
$ws.Range("N110").Value = -8537.5
$ws.Range("H136").Value = 4167385
$ws.Range("I136").Value = 4167385
$ws.Range("K136").Value = 12502155
$ws.Range("M136").Value = -12499605

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1647
$ws.Range("I20").Value = 1337.3636
$ws.Range("K20").Value = 1337.3636
$ws.Range("M20").Value = -1090.3636
$ws.Range("H80").Value = 5519.125
$ws.Range("I80").Value = 17283
$ws.Range("J80").Value = 2804.3845
$ws.Range("K80").Value = 17283
$ws.Range("L80").Value = 2804.3845
$ws.Range("M80").Value = -16285
$ws.Range("N80").Value = -4800.3845
$ws.Range("H83").Value = 5519.125
$ws.Range("I83").Value = 17283
$ws.Range("J83").Value = 2804.3845
$ws.Range("K83").Value = 86415
$ws.Range("L83").Value = 14021.9225
$ws.Range("M83").Value = -81423
$ws.Range("N83").Value = -24005.9225
$ws.Range("H107").Value = 4188.478
$ws.Range("I107").Value = 3576.2144
$ws.Range("K107").Value = 3576.2144
$ws.Range("M107").Value = -1656.2144
$ws.Range("H134").Value = 730405.25
$ws.Range("I134").Value = 796105.0600000001
$ws.Range("J134").Value = 511406
$ws.Range("K134").Value = 2388315.18
$ws.Range("L134").Value = 1534218
$ws.Range("M134").Value = -2385780.18
$ws.Range("N134").Value = -1539288

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10622.1455
$ws.Range("I31").Value = 3757.0881
$ws.Range("K31").Value = 3757.0881
$ws.Range("M31").Value = -3462.0881
$ws.Range("H34").Value = 10622.1455
$ws.Range("I34").Value = 3757.0881
$ws.Range("K34").Value = 3757.0881
$ws.Range("M34").Value = -3555.0881
$ws.Range("H58").Value = 1380484.9
$ws.Range("I58").Value = 2472873.8
$ws.Range("J58").Value = 14998.75
$ws.Range("K58").Value = 2472873.8
$ws.Range("L58").Value = 14998.75
$ws.Range("M58").Value = -2472670.8
$ws.Range("N58").Value = -15404.75
$ws.Range("H82").Value = 63725
$ws.Range("J82").Value = 63725
$ws.Range("L82").Value = 63725
$ws.Range("N82").Value = -64447
$ws.Range("H85").Value = 63725
$ws.Range("J85").Value = 63725
$ws.Range("L85").Value = 63725
$ws.Range("N85").Value = -66221
$ws.Range("H86").Value = 93232.87
$ws.Range("I86").Value = 5890.75
$ws.Range("J86").Value = 188515.19
$ws.Range("K86").Value = 5890.75
$ws.Range("L86").Value = 188515.19
$ws.Range("M86").Value = -4767.75
$ws.Range("N86").Value = -190761.19
$ws.Range("H89").Value = 93232.87
$ws.Range("I89").Value = 5890.75
$ws.Range("J89").Value = 188515.19
$ws.Range("K89").Value = 29453.75
$ws.Range("L89").Value = 942575.95
$ws.Range("M89").Value = -23837.75
$ws.Range("N89").Value = -953807.95
$ws.Range("H99").Value = 3341.611
$ws.Range("I99").Value = 2870
$ws.Range("J99").Value = 3813.2222
$ws.Range("K99").Value = 2870
$ws.Range("L99").Value = 3813.2222
$ws.Range("M99").Value = -1372
$ws.Range("N99").Value = -6809.2222
$ws.Range("H105").Value = 33304.547
$ws.Range("I105").Value = 44669
$ws.Range("K105").Value = 44669
$ws.Range("M105").Value = -42922
$ws.Range("H126").Value = 3341.611
$ws.Range("I126").Value = 2870
$ws.Range("J126").Value = 3813.2222
$ws.Range("K126").Value = 8610
$ws.Range("L126").Value = 11439.6666
$ws.Range("M126").Value = -6140
$ws.Range("N126").Value = -16379.6666
$ws.Range("H132").Value = 41883700
$ws.Range("I132").Value = 58826124
$ws.Range("K132").Value = 176478372
$ws.Range("M132").Value = -176475842
$ws.Range("H134").Value = 6286.3447
$ws.Range("I134").Value = 6743.846
$ws.Range("J134").Value = 2321.3333
$ws.Range("K134").Value = 20231.538
$ws.Range("L134").Value = 6963.999899999999
$ws.Range("M134").Value = -17696.538
$ws.Range("N134").Value = -12033.9999
$ws.Range("H136").Value = 1380484.9
$ws.Range("I136").Value = 2472873.8
$ws.Range("J136").Value = 14998.75
$ws.Range("K136").Value = 7418621.399999999
$ws.Range("L136").Value = 44996.25
$ws.Range("M136").Value = -7416071.399999999
$ws.Range("N136").Value = -50096.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 39
$ws.Range("I12").Value = 23
$ws.Range("K12").Value = 69
$ws.Range("M12").Value = 104
$ws.Range("H26").Value = 472.55554
$ws.Range("J26").Value = 447
$ws.Range("L26").Value = 1341
$ws.Range("N26").Value = -1917
$ws.Range("H38").Value = 34.75
$ws.Range("I38").Value = 33
$ws.Range("K38").Value = 99
$ws.Range("M38").Value = 248

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4968
$ws.Range("I70").Value = 4964.4443
$ws.Range("K70").Value = 4964.4443
$ws.Range("M70").Value = -4694.4443
$ws.Range("H73").Value = 4968
$ws.Range("I73").Value = 4964.4443
$ws.Range("K73").Value = 4964.4443
$ws.Range("M73").Value = -4028.4443
$ws.Range("H97").Value = 866.7778
$ws.Range("J97").Value = 643.2
$ws.Range("L97").Value = 643.2
$ws.Range("N97").Value = -1635.2
$ws.Range("H132").Value = 29771398
$ws.Range("I132").Value = 40486090
$ws.Range("J132").Value = 8367.777
$ws.Range("K132").Value = 121458270
$ws.Range("L132").Value = 25103.331
$ws.Range("M132").Value = -121455740
$ws.Range("N132").Value = -30163.331
$ws.Range("H135").Value = 94918.664
$ws.Range("J135").Value = 94918.664
$ws.Range("L135").Value = 94918.664
$ws.Range("N135").Value = -105058.664

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4233.8423
$ws.Range("I40").Value = 4191.278
$ws.Range("K40").Value = 4191.278
$ws.Range("M40").Value = -4055.278
$ws.Range("H47").Value = 30000
$ws.Range("J47").Value = 30000
$ws.Range("L47").Value = 30000
$ws.Range("N47").Value = -30980
$ws.Range("H52").Value = 30000
$ws.Range("J52").Value = 30000
$ws.Range("L52").Value = 30000
$ws.Range("N52").Value = -30466
$ws.Range("H136").Value = 3195.611
$ws.Range("I136").Value = 1956.909
$ws.Range("J136").Value = 5142.143
$ws.Range("K136").Value = 5870.727000000001
$ws.Range("L136").Value = 15426.429
$ws.Range("M136").Value = -3320.727000000001
$ws.Range("N136").Value = -20526.429

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 1000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H132").Value = 9585927
$ws.Range("I132").Value = 10593498
$ws.Range("K132").Value = 31780494
$ws.Range("M132").Value = -31777964
